$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5288504903118678
$ws.Cells.Item(2, 3).Value = 0.09332851332715109
$ws.Cells.Item(2, 5).Value = 0.0978206075889716
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.002456794523338966
$ws.Cells.Item(2, 9).Value = 0.7969589987911156
$ws.Cells.Item(2, 11).Value = 0.3078994974998182
$ws.Cells.Item(2, 12).Value = 0.2051546940440119
$ws.Cells.Item(2, 14).Value = 1.642873267952641
$ws.Cells.Item(2, 15).Value = 3.05790235997209

$ws.Cells.Item(3, 2).Value = 0.4907794919976425
$ws.Cells.Item(3, 3).Value = 0.09249425262155597
$ws.Cells.Item(3, 5).Value = 0.09712084892509409
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.002458877658145564
$ws.Cells.Item(3, 9).Value = 0.804172606307116
$ws.Cells.Item(3, 11).Value = 0.2748926563203611
$ws.Cells.Item(3, 12).Value = 0.1980262487047781
$ws.Cells.Item(3, 14).Value = 1.658551352523155
$ws.Cells.Item(3, 15).Value = 3.081946286615079

$ws.Cells.Item(4, 2).Value = 0.4675359107358474
$ws.Cells.Item(4, 3).Value = 0.09197939975822322
$ws.Cells.Item(4, 5).Value = 0.09674053665384719
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.002460225474475847
$ws.Cells.Item(4, 9).Value = 0.8089730015781953
$ws.Cells.Item(4, 11).Value = 0.2546364668096004
$ws.Cells.Item(4, 12).Value = 0.1937564448849542
$ws.Cells.Item(4, 14).Value = 1.668707977147392
$ws.Cells.Item(4, 15).Value = 3.09825133687724

$ws.Cells.Item(5, 2).Value = 0.4580978944734397
$ws.Cells.Item(5, 3).Value = 0.09176895142292096
$ws.Cells.Item(5, 5).Value = 0.09659799159290827
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.002460792060917376
$ws.Cells.Item(5, 9).Value = 0.8110225748872679
$ws.Cells.Item(5, 11).Value = 0.2463849804072993
$ws.Cells.Item(5, 12).Value = 0.1920434963982558
$ws.Cells.Item(5, 14).Value = 1.672980275585793
$ws.Cells.Item(5, 15).Value = 3.105283602864219

$ws.Cells.Item(6, 2).Value = 0.4565327889958155
$ws.Cells.Item(6, 3).Value = 0.09173396834481196
$ws.Cells.Item(6, 5).Value = 0.09657507396922682
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.002460887190743464
$ws.Cells.Item(6, 9).Value = 0.8113685453529484
$ws.Cells.Item(6, 11).Value = 0.245015028552956
$ws.Cells.Item(6, 12).Value = 0.191760698239861
$ws.Cells.Item(6, 14).Value = 1.67369774592067
$ws.Cells.Item(6, 15).Value = 3.106474731415318

$ws.Cells.Item(7, 2).Value = 0.4674084880348346
$ws.Cells.Item(7, 3).Value = 0.09197656415561539
$ws.Cells.Item(7, 5).Value = 0.09673856385434831
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.002460233045360803
$ws.Cells.Item(7, 9).Value = 0.8090002647369232
$ws.Cells.Item(7, 11).Value = 0.2545251712313217
$ws.Cells.Item(7, 12).Value = 0.1937332338780919
$ws.Cells.Item(7, 14).Value = 1.668765054675493
$ws.Cells.Item(7, 15).Value = 3.098344606075756

$ws.Cells.Item(8, 2).Value = 0.51569658531497
$ws.Cells.Item(8, 3).Value = 0.09304141171787705
$ws.Cells.Item(8, 5).Value = 0.09756910598700586
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.002457498546505719
$ws.Cells.Item(8, 9).Value = 0.7993692234903591
$ws.Cells.Item(8, 11).Value = 0.2965169376196855
$ws.Cells.Item(8, 12).Value = 0.2026746383695865
$ws.Cells.Item(8, 14).Value = 1.648169064462198
$ws.Cells.Item(8, 15).Value = 3.065872715522005

$ws.Cells.Item(9, 2).Value = 0.611412287532346
$ws.Cells.Item(9, 3).Value = 0.09510825153324021
$ws.Cells.Item(9, 5).Value = 0.09958828279118848
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.002452679531530162
$ws.Cells.Item(9, 9).Value = 0.7834269097077815
$ws.Cells.Item(9, 11).Value = 0.3789241049839234
$ws.Cells.Item(9, 12).Value = 0.2210550141449801
$ws.Cells.Item(9, 14).Value = 1.6119840819958
$ws.Cells.Item(9, 15).Value = 3.01442759497165

$ws.Cells.Item(10, 2).Value = 0.6823307112372561
$ws.Cells.Item(10, 3).Value = 0.09661311320139987
$ws.Cells.Item(10, 5).Value = 0.1013087351629736
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.002449467038568741
$ws.Cells.Item(10, 9).Value = 0.773507241100571
$ws.Cells.Item(10, 11).Value = 0.4394863795191952
$ws.Cells.Item(10, 12).Value = 0.235072049503998
$ws.Cells.Item(10, 14).Value = 1.587956372316899
$ws.Cells.Item(10, 15).Value = 2.98408538986088

$ws.Cells.Item(11, 2).Value = 0.7147173687435782
$ws.Cells.Item(11, 3).Value = 0.09729462296643732
$ws.Cells.Item(11, 5).Value = 0.1021426537420353
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.002448076143938867
$ws.Cells.Item(11, 9).Value = 0.7693835557006139
$ws.Cells.Item(11, 11).Value = 0.4670377363888463
$ws.Cells.Item(11, 12).Value = 0.2415596120729333
$ws.Cells.Item(11, 14).Value = 1.57757974998626
$ws.Cells.Item(11, 15).Value = 2.971900212845725

$ws.Cells.Item(12, 2).Value = 0.7269987836929772
$ws.Cells.Item(12, 3).Value = 0.09755223921934686
$ws.Cells.Item(12, 5).Value = 0.1024657863305158
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.002447559532917669
$ws.Cells.Item(12, 9).Value = 0.7678779139055294
$ws.Cells.Item(12, 11).Value = 0.4774704428247958
$ws.Cells.Item(12, 12).Value = 0.2440321825285281
$ws.Cells.Item(12, 14).Value = 1.57372995874621
$ws.Cells.Item(12, 15).Value = 2.967518585337217

$ws.Cells.Item(13, 2).Value = 0.7243530043013777
$ws.Cells.Item(13, 3).Value = 0.09749677752773778
$ws.Cells.Item(13, 5).Value = 0.1023958676862939
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.002447670346410573
$ws.Cells.Item(13, 9).Value = 0.7681996943237657
$ws.Cells.Item(13, 11).Value = 0.475223599572729
$ws.Cells.Item(13, 12).Value = 0.2434989661332025
$ws.Cells.Item(13, 14).Value = 1.574555539631191
$ws.Cells.Item(13, 15).Value = 2.968451900428306

$ws.Cells.Item(14, 2).Value = 0.7157274259532755
$ws.Cells.Item(14, 3).Value = 0.09731582646018921
$ws.Cells.Item(14, 5).Value = 0.1021690909947175
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.002448033440278111
$ws.Cells.Item(14, 9).Value = 0.769258565177811
$ws.Cells.Item(14, 11).Value = 0.467896052641521
$ws.Cells.Item(14, 12).Value = 0.2417627144851053
$ws.Cells.Item(14, 14).Value = 1.5772614298889
$ws.Cells.Item(14, 15).Value = 2.971535071260291

$ws.Cells.Item(15, 2).Value = 0.710446239418161
$ws.Cells.Item(15, 3).Value = 0.09720492873310604
$ws.Cells.Item(15, 5).Value = 0.1020311395276288
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.002448257158196343
$ws.Cells.Item(15, 9).Value = 0.7699144349620823
$ws.Cells.Item(15, 11).Value = 0.4634076500474293
$ws.Cells.Item(15, 12).Value = 0.2407012738274972
$ws.Cells.Item(15, 14).Value = 1.578929233135035
$ws.Cells.Item(15, 15).Value = 2.973453898877921

$ws.Cells.Item(16, 2).Value = 0.6802166278376944
$ws.Cells.Item(16, 3).Value = 0.09656851193842186
$ws.Cells.Item(16, 5).Value = 0.1012552660629105
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.002449559351159193
$ws.Cells.Item(16, 9).Value = 0.7737845574052891
$ws.Cells.Item(16, 11).Value = 0.4376858098140985
$ws.Cells.Item(16, 12).Value = 0.2346502990472459
$ws.Cells.Item(16, 14).Value = 1.588645649880888
$ws.Cells.Item(16, 15).Value = 2.984914272627123

$ws.Cells.Item(17, 2).Value = 0.6617033475613709
$ws.Cells.Item(17, 3).Value = 0.09617729531596098
$ws.Cells.Item(17, 5).Value = 0.1007924061922374
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.002450376224913065
$ws.Cells.Item(17, 9).Value = 0.7762583360427584
$ws.Cells.Item(17, 11).Value = 0.421906209348947
$ws.Cells.Item(17, 12).Value = 0.2309666096744536
$ws.Cells.Item(17, 14).Value = 1.594748164604965
$ws.Cells.Item(17, 15).Value = 2.992359172003802

$ws.Cells.Item(18, 2).Value = 0.6510668573799308
$ws.Cells.Item(18, 3).Value = 0.09595199089859108
$ws.Cells.Item(18, 5).Value = 0.1005310094590826
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.002450852706100071
$ws.Cells.Item(18, 9).Value = 0.7777177860292781
$ws.Cells.Item(18, 11).Value = 0.4128303550119483
$ws.Cells.Item(18, 12).Value = 0.2288583169445531
$ws.Cells.Item(18, 14).Value = 1.598310287363688
$ws.Cells.Item(18, 15).Value = 2.996793540400034

$ws.Cells.Item(19, 2).Value = 0.6474675844392834
$ws.Cells.Item(19, 3).Value = 0.09587565792092079
$ws.Cells.Item(19, 5).Value = 0.1004433350767862
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.002451015175709653
$ws.Cells.Item(19, 9).Value = 0.7782182165652038
$ws.Cells.Item(19, 11).Value = 0.4097574708147818
$ws.Cells.Item(19, 12).Value = 0.2281462865030335
$ws.Cells.Item(19, 14).Value = 1.599525312580216
$ws.Cells.Item(19, 15).Value = 2.998321091421033

$ws.Cells.Item(20, 2).Value = 0.6636728959795732
$ws.Cells.Item(20, 3).Value = 0.09621897081500208
$ws.Cells.Item(20, 5).Value = 0.1008411789340471
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.002450288580845115
$ws.Cells.Item(20, 9).Value = 0.7759912103608571
$ws.Cells.Item(20, 11).Value = 0.4235859629508241
$ws.Cells.Item(20, 12).Value = 0.231357662102937
$ws.Cells.Item(20, 14).Value = 1.594093147598844
$ws.Cells.Item(20, 15).Value = 2.991550891359921

$ws.Cells.Item(21, 2).Value = 0.7182605041779766
$ws.Cells.Item(21, 3).Value = 0.0973689887034439
$ws.Cells.Item(21, 5).Value = 0.1022355016670673
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.002447926517119837
$ws.Cells.Item(21, 9).Value = 0.7689460318325416
$ws.Cells.Item(21, 11).Value = 0.4700483461082854
$ws.Cells.Item(21, 12).Value = 0.242272263826294
$ws.Cells.Item(21, 14).Value = 1.576464483802752
$ws.Cells.Item(21, 15).Value = 2.970623155461396

$ws.Cells.Item(22, 2).Value = 0.7540370930472022
$ws.Cells.Item(22, 3).Value = 0.09811792211020531
$ws.Cells.Item(22, 5).Value = 0.1031895756938823
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.002446441567024211
$ws.Cells.Item(22, 9).Value = 0.764667475515278
$ws.Cells.Item(22, 11).Value = 0.5004116319315699
$ws.Cells.Item(22, 12).Value = 0.2494980451919133
$ws.Cells.Item(22, 14).Value = 1.565407159950603
$ws.Cells.Item(22, 15).Value = 2.958301651418452

$ws.Cells.Item(23, 2).Value = 0.7349335151694447
$ws.Cells.Item(23, 3).Value = 0.09771845206143581
$ws.Cells.Item(23, 5).Value = 0.1026764611026856
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.002447228747823461
$ws.Cells.Item(23, 9).Value = 0.7669212051467156
$ws.Cells.Item(23, 11).Value = 0.4842066006997641
$ws.Cells.Item(23, 12).Value = 0.2456330865761629
$ws.Cells.Item(23, 14).Value = 1.571266206409458
$ws.Cells.Item(23, 15).Value = 2.964753791320874

$ws.Cells.Item(24, 2).Value = 0.6627824411492611
$ws.Cells.Item(24, 3).Value = 0.0962001305312512
$ws.Cells.Item(24, 5).Value = 0.1008191141121166
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.002450328183371983
$ws.Cells.Item(24, 9).Value = 0.776111861876231
$ws.Cells.Item(24, 11).Value = 0.4228265586134228
$ws.Cells.Item(24, 12).Value = 0.2311808376722269
$ws.Cells.Item(24, 14).Value = 1.594389113495556
$ws.Cells.Item(24, 15).Value = 2.991915834728943

$ws.Cells.Item(25, 2).Value = 0.5854119878206063
$ws.Cells.Item(25, 3).Value = 0.0945514664821232
$ws.Cells.Item(25, 5).Value = 0.09900034973895799
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.002453925366833519
$ws.Cells.Item(25, 9).Value = 0.7874247123259011
$ws.Cells.Item(25, 11).Value = 0.3566262444002746
$ws.Cells.Item(25, 12).Value = 0.2159923360058542
$ws.Cells.Item(25, 14).Value = 1.621323601154312
$ws.Cells.Item(25, 15).Value = 3.027035589446626
